# Swap the contents of data rows 2 and 3 on the active sheet: row 2 becomes
# the record that used to be row 3 (Nordfladdermus / Eptesicus nilssonii),
# and row 3 becomes the record that used to be row 2 (Violmussling /
# Trichaptum laricinum). Cells that already hold identical content in both
# rows (C, D, P, S, T, U, V, W, Z, AB, AD, AE, AG, AT, AY) are left alone.

function Set-TextCell($range, $text) {
    # Force a genuine text value even when it looks like a number or date,
    # then restore the default (un-styled) look of the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (becomes the old row 3 record) ----
$ws.Range("A2").Value = 63191078
$ws.Range("B2").Value = 57484
$ws.Range("E2").Value = 205998
Set-TextCell $ws.Range("F2") "Nordfladdermus"
Set-TextCell $ws.Range("G2") "Eptesicus nilssonii"
Set-TextCell $ws.Range("H2") "(A.Keyserling & Blasius, 1839)"
Set-TextCell $ws.Range("I2") "9"
Set-TextCell $ws.Range("J2") "registreringar"
Set-TextCell $ws.Range("K2") "adult"
Set-TextCell $ws.Range("M2") "födosökande"
Set-TextCell $ws.Range("N2") "autobox med höghastighetsinspelning"
$ws.Range("Q2").Value = 759836.9302479513
$ws.Range("R2").Value = 7086399.103862511
Set-TextCell $ws.Range("Y2") "2016-06-28"
Set-TextCell $ws.Range("AA2") "2016-07-01"
Set-TextCell $ws.Range("AC2") "lokal-ID: 5"
Set-TextCell $ws.Range("AW2") "Alexander Eriksson"
Set-TextCell $ws.Range("AX2") "Alexander Eriksson, Björn Palmqvist"

# ---- Row 3 (becomes the old row 2 record) ----
$ws.Range("A3").Value = 69062005
$ws.Range("B3").Value = 89557
$ws.Range("E3").Value = 1588
Set-TextCell $ws.Range("F3") "Violmussling"
Set-TextCell $ws.Range("G3") "Trichaptum laricinum"
Set-TextCell $ws.Range("H3") "(P.Karst.) Ryvarden"
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("Q3").Value = 759914.2062387039
$ws.Range("R3").Value = 7086527.972586548
Set-TextCell $ws.Range("Y3") "2017-08-29"
Set-TextCell $ws.Range("AA3") "2017-08-29"
$ws.Range("AC3").ClearContents()
Set-TextCell $ws.Range("AW3") "Torbjörn Josefsson"
Set-TextCell $ws.Range("AX3") "Torbjörn Josefsson, Karin Björnehall"
